# Updated cryptos list with GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Force the cell to be stored as text (matches the source data's
    # inline-string Price column) instead of letting Excel auto-coerce
    # numeric-looking strings ("235.17", "59.00", ...) into numbers.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "37.389.60"
$ws.Range("E2").Value = "  +2.17%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.066.77"
$ws.Range("E3").Value = "  +3.35%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - BNB
Set-TextCell "D5" "235.17"
$ws.Range("E5").Value = "  +0.36%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +2.15%  "

# Row 7 - Solana
Set-TextCell "D7" "58.28"
$ws.Range("E7").Value = "  +5.62%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.382"
$ws.Range("E9").Value = "  +2.87%  "

# Row 10 - OKB
Set-TextCell "D10" "59.00"
$ws.Range("E10").Value = "  +1.52%  "

# Row 11 - Dogecoin
Set-TextCell "D11" "0.0763"
$ws.Range("E11").Value = "  +1.86%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +2.78%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextCell "D13" "2.369.97"
$ws.Range("E13").Value = "  +3.46%  "

# Row 14 - Chainlink
$ws.Range("E14").Value = "  +2.15%  "

# Row 15 - Avalanche
Set-TextCell "D15" "21.31"
$ws.Range("E15").Value = "  +4.14%  "

# Row 16 - Polygon
$ws.Range("E16").Value = "  +2.33%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  +2.06%  "

# Row 18 - WrappedEther
Set-TextCell "D18" "2.025.95"
$ws.Range("E18").Value = "  +1.31%  "

# Row 19 - WrappedBTC
Set-TextCell "D19" "37.365.52"
$ws.Range("E19").Value = "  +2.33%  "

# Row 20 - Uniswap
Set-TextCell "D20" "6.16"
$ws.Range("E20").Value = "  +15.51%  "

# Row 21 - Litecoin
Set-TextCell "D21" "70.16"
$ws.Range("E21").Value = "  +3.43%  "

# Row 22 - ShibaInu
$ws.Range("E22").Value = "  +1.00%  "

# Row 23 - BitcoinCash
Set-TextCell "D23" "227.26"
$ws.Range("E23").Value = "  +2.23%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.00%  "

# Row 25 - PancakeSwap
Set-TextCell "D25" "2.44"
$ws.Range("E25").Value = "  +1.77%  "

# Row 26 - Toncoin
Set-TextCell "D26" "2.39"
$ws.Range("E26").Value = "  +0.78%  "

# Row 27 - Monero
Set-TextCell "D27" "165.31"
$ws.Range("E27").Value = "  +2.03%  "

# Row 28 - ImmutableX
$ws.Range("E28").Value = "  +13.79%  "

# Row 29 - Cosmos
Set-TextCell "D29" "8.88"
$ws.Range("E29").Value = "  +2.23%  "

# Row 30 - EthereumClassic
$ws.Range("E30").Value = "  +2.09%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  -0.60%  "

# Row 32 - Stellar
$ws.Range("E32").Value = "  +1.93%  "

# Row 33 - Filecoin
Set-TextCell "D33" "4.54"
$ws.Range("E33").Value = "  +3.77%  "

# Row 34 - Hedera
Set-TextCell "D34" "0.0622"
$ws.Range("E34").Value = "  +2.82%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +8.23%  "

# Row 36 - InternetComputer(DFINITY)
$ws.Range("E36").Value = "  +7.07%  "

# Row 37 - RenderToken
Set-TextCell "D37" "3.37"
$ws.Range("E37").Value = "  -0.02%  "

# Row 38 - BinanceUSD
$ws.Range("E38").Value = "  +0.05%  "

# Row 39 - WEMIXToken
$ws.Range("E39").Value = "  +1.20%  "

# Row 40 - THORChain
$ws.Range("E40").Value = "  +3.35%  "

# Row 41 - Cronos
Set-TextCell "D41" "0.0973"
$ws.Range("E41").Value = "  +2.71%  "

# Row 42 - HuobiToken
$ws.Range("E42").Value = "  -1.41%  "

# Row 43 - FTXToken
$ws.Range("E43").Value = "  +22.65%  "

# Row 44 - Maker
Set-TextCell "D44" "1.459.20"
$ws.Range("E44").Value = "  +0.20%  "

# Row 45 - Aave
Set-TextCell "D45" "95.83"
$ws.Range("E45").Value = "  +7.32%  "

# Row 46 - now TrustWalletToken (was VeChain)
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D46" "1.16"
$ws.Range("E46").Value = "  +5.89%  "

# Row 47 - now VeChain (was TrustWalletToken)
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D47" "0.0211"
$ws.Range("E47").Value = "  +4.03%  "

# Row 48 - InjectiveProtocol
Set-TextCell "D48" "15.84"
$ws.Range("E48").Value = "  +3.79%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  +3.91%  "

# Row 50 - FraxShare
Set-TextCell "D50" "7.29"
$ws.Range("E50").Value = "  +6.49%  "

# Row 51 - MXToken
$ws.Range("E51").Value = "  +2.08%  "
